# Rename "Interventions coverages" -> "Interventions cost and coverage"
$wb = $excel.ActiveWorkbook

$covSheet = $wb.Worksheets.Item("Interventions coverages")
$covSheet.Name = "Interventions cost and coverage"

# Delete the birth-time/order related worksheets
$namesToDelete = @("RR birth by type", "birth distribution", "time between births", "RR birth by time")
foreach ($name in $namesToDelete) {
    $sheet = $wb.Worksheets.Item($name)
    $sheet.Delete()
}
